$d = $word.ActiveDocument

# 1. Update the date heading
$d.Content.Find.Execute("May 07, 2024", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "May 13, 2024", 2) | Out-Null

# 2. Update the italic instructions paragraph (also drops the stray trailing space)
$d.Content.Find.Execute(
    "Extract any quote mentions " + [char]0x201C + "{variable_name}" + [char]0x201D + ". Only include direct quotation with the corresponding page number(s). ",
    $true, $false, $false, $false, $false, `
    $true, 1, $false,
    "Extract any quote that addresses " + [char]0x201C + "{variable_name}" + [char]0x201D + ". Only include direct quotations with the corresponding page number(s).",
    2) | Out-Null

# 3. Variable name cell in the first (summary) table: "Cement" -> "energy"
$summaryTable = $d.Tables.Item(1)
$summaryTable.Cell(2, 1).Range.Text = "energy"

# 4. Source document heading: IRE03 ... -> Argentina ...
$d.Content.Find.Execute("IRE03 CLIMATE ACTION PLAN 2023.pdf", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Argentina_Actualización meta de emisiones 2030.pdf", 2) | Out-Null

# 5. Quote table: "Related Variables" first column cell "Cement" -> "energy"
$quoteTable = $d.Tables.Item(2)
$quoteTable.Cell(2, 1).Range.Text = "energy"

# 6. Replace the big Quote cell contents with the new (Spanish) quotes
$quoteCellRange = $quoteTable.Cell(2, 2).Range
$start = $quoteCellRange.Start
$end = $quoteCellRange.End
$newQuoteText = "En el caso de las variables relacionadas con el sector energético, se utilizaron los mismos modelos de demanda y oferta de energía utilizados para la planificación energética nacional. [page 12]. " + `
    [char]11 + `
    "Se contempló una demanda creciente del consumo por parte de la población, con medidas de eficiencia energética en todos los sectores, un aumento significativo del porcentaje de energías renovables y de generación distribuida, y una mayor producción de gas natural en términos absolutos y relativos respecto a la producción de petróleo. [page 12]. " + `
    [char]11 + `
    "Entre ellas, cabe mencionar la población, el producto interno bruto, la demanda y la oferta de energía, las existencias de ganado bovino, la producción agrícola y el cambio de uso de la tierra. [page 12]. " + `
    [char]11 + `
    "En cuenta políticas activas tendientes a aumentar la eficiencia y la utilización de gas natural y electricidad. [page 12]. " + `
    [char]11
$d.Range($start, $end).Text = $newQuoteText

# 7. Footer stats paragraph
$d.Content.Find.Execute("1 documents (284 total pages) processed in 55.10 seconds", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1 documents (18 total pages) processed in 23.66 seconds", 2) | Out-Null
